# Build site at 2023-04-12 14:53:07 UTC
# Rewrite LOQ4230 "Estagio em Engenharia de Producao" course description block
# (rows 13-42 of sheet1) to the updated copy, inserting a new row and
# re-flowing the remaining rows; rows 1-12 are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: make sure every cell we are about to (re)write in rows 13-42 carries
# the correct existing style (A=1 bold label, B=2 wrapped text, C=3 red wrapped
# text) by copying format from the permanent, untouched donor cells in row 3.
$ws.Range("B3").Copy($ws.Range("B13")) | Out-Null
$ws.Range("C3").Copy($ws.Range("C13")) | Out-Null
$ws.Range("A3").Copy($ws.Range("A14")) | Out-Null
$ws.Range("B3").Copy($ws.Range("B14")) | Out-Null
$ws.Range("C3").Copy($ws.Range("C14")) | Out-Null
$ws.Range("A3").Copy($ws.Range("A15")) | Out-Null
$ws.Range("B3").Copy($ws.Range("B15")) | Out-Null
$ws.Range("C3").Copy($ws.Range("C15")) | Out-Null
$ws.Range("A3").Copy($ws.Range("A16")) | Out-Null
$ws.Range("B3").Copy($ws.Range("B16")) | Out-Null
$ws.Range("C3").Copy($ws.Range("C16")) | Out-Null
$ws.Range("A3").Copy($ws.Range("A17")) | Out-Null
$ws.Range("B3").Copy($ws.Range("B17")) | Out-Null
$ws.Range("C3").Copy($ws.Range("C17")) | Out-Null
$ws.Range("A3").Copy($ws.Range("A18")) | Out-Null
$ws.Range("A3").Copy($ws.Range("A19")) | Out-Null
$ws.Range("B3").Copy($ws.Range("B19")) | Out-Null
$ws.Range("C3").Copy($ws.Range("C19")) | Out-Null
$ws.Range("A3").Copy($ws.Range("A20")) | Out-Null
$ws.Range("B3").Copy($ws.Range("B20")) | Out-Null
$ws.Range("C3").Copy($ws.Range("C20")) | Out-Null
$ws.Range("A3").Copy($ws.Range("A21")) | Out-Null
$ws.Range("B3").Copy($ws.Range("B21")) | Out-Null
$ws.Range("C3").Copy($ws.Range("C21")) | Out-Null
$ws.Range("A3").Copy($ws.Range("A22")) | Out-Null
$ws.Range("B3").Copy($ws.Range("B22")) | Out-Null
$ws.Range("C3").Copy($ws.Range("C22")) | Out-Null
$ws.Range("A3").Copy($ws.Range("A23")) | Out-Null
$ws.Range("B3").Copy($ws.Range("B24")) | Out-Null
$ws.Range("C3").Copy($ws.Range("C24")) | Out-Null
$ws.Range("B3").Copy($ws.Range("B25")) | Out-Null
$ws.Range("C3").Copy($ws.Range("C25")) | Out-Null
$ws.Range("B3").Copy($ws.Range("B26")) | Out-Null
$ws.Range("C3").Copy($ws.Range("C26")) | Out-Null
$ws.Range("B3").Copy($ws.Range("B27")) | Out-Null
$ws.Range("C3").Copy($ws.Range("C27")) | Out-Null
$ws.Range("B3").Copy($ws.Range("B28")) | Out-Null
$ws.Range("C3").Copy($ws.Range("C28")) | Out-Null
$ws.Range("B3").Copy($ws.Range("B29")) | Out-Null
$ws.Range("C3").Copy($ws.Range("C29")) | Out-Null
$ws.Range("B3").Copy($ws.Range("B30")) | Out-Null
$ws.Range("C3").Copy($ws.Range("C30")) | Out-Null
$ws.Range("B3").Copy($ws.Range("B31")) | Out-Null
$ws.Range("C3").Copy($ws.Range("C31")) | Out-Null
$ws.Range("B3").Copy($ws.Range("B32")) | Out-Null
$ws.Range("C3").Copy($ws.Range("C32")) | Out-Null
$ws.Range("B3").Copy($ws.Range("B33")) | Out-Null
$ws.Range("C3").Copy($ws.Range("C33")) | Out-Null
$ws.Range("B3").Copy($ws.Range("B34")) | Out-Null
$ws.Range("C3").Copy($ws.Range("C34")) | Out-Null
$ws.Range("B3").Copy($ws.Range("B35")) | Out-Null
$ws.Range("C3").Copy($ws.Range("C35")) | Out-Null
$ws.Range("B3").Copy($ws.Range("B36")) | Out-Null
$ws.Range("C3").Copy($ws.Range("C36")) | Out-Null
$ws.Range("B3").Copy($ws.Range("B37")) | Out-Null
$ws.Range("C3").Copy($ws.Range("C37")) | Out-Null
$ws.Range("B3").Copy($ws.Range("B38")) | Out-Null
$ws.Range("C3").Copy($ws.Range("C38")) | Out-Null
$ws.Range("B3").Copy($ws.Range("B39")) | Out-Null
$ws.Range("C3").Copy($ws.Range("C39")) | Out-Null
$ws.Range("B3").Copy($ws.Range("B40")) | Out-Null
$ws.Range("C3").Copy($ws.Range("C40")) | Out-Null
$ws.Range("B3").Copy($ws.Range("B41")) | Out-Null
$ws.Range("C3").Copy($ws.Range("C41")) | Out-Null
$ws.Range("B3").Copy($ws.Range("B42")) | Out-Null
$ws.Range("C3").Copy($ws.Range("C42")) | Out-Null

# Step 2: write the new text content for rows 13-42
$ws.Range("B13").Value = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Range("C13").Value = "5840560 - Marco Antonio Carvalho Pereira"

$ws.Range("A14").Value = "Programa resumido:"
$ws.Range("B14").Value = "Plano de Trabalho específico. Realização do Estágio. Relatório final e/ou parciais."
$ws.Range("C14").Value = "Plano de Trabalho específico. Realização do Estágio. Relatório final e/ou parciais."

$ws.Range("A15").Value = "Short syllabus:"
$ws.Range("B15").Value = "Specific Work Plan. Realization of the Internship. Final and / or partial report."
$ws.Range("C15").Value = "Specific Work Plan. Realization of the Internship. Final and / or partial report."

$ws.Range("A16").Value = "Programa:"
$ws.Range("B16").Value = "Participação do aluno em processo seletivo de empresas ou no setor acadêmico. Estágio realizado sob a supervisão da Escola de Engenharia de Lorena, através do Departamento em Engenharia Química. O conteúdo será estabelecido individualmente no Plano de Trabalho entre o Supervisor do Estágio e o professor orientador, desde que relacionado com as áreas afins da Engenharia de Produção. Apresentação de relatório final e/ou relatórios parciais sobre as atividades desenvolvidas no estágio."
$ws.Range("C16").Value = "Participação do aluno em processo seletivo de empresas ou no setor acadêmico. Estágio realizado sob a supervisão da Escola de Engenharia de Lorena, através do Departamento em Engenharia Química. O conteúdo será estabelecido individualmente no Plano de Trabalho entre o Supervisor do Estágio e o professor orientador, desde que relacionado com as áreas afins da Engenharia de Produção. Apresentação de relatório final e/ou relatórios parciais sobre as atividades desenvolvidas no estágio."

$ws.Range("A17").Value = "Syllabus:"
$ws.Range("B17").Value = "Participation of the student in the selective process of companies or in the academic sector. Internship carried out under the supervision of the School of Engineering of Lorena, through the Department of Chemical Engineering. The content will be established individually in the Work Plan between the Internship Supervisor and the tutor, as long as related to the areas of Industrial Engineering. Presentation of final report and / or partial reports about the activities carried out during the internship."
$ws.Range("C17").Value = "Participation of the student in the selective process of companies or in the academic sector. Internship carried out under the supervision of the School of Engineering of Lorena, through the Department of Chemical Engineering. The content will be established individually in the Work Plan between the Internship Supervisor and the tutor, as long as related to the areas of Industrial Engineering. Presentation of final report and / or partial reports about the activities carried out during the internship."

$ws.Range("A18").Value = "Avaliação:"

$ws.Range("A19").Value = "Método:"
$ws.Range("B19").Value = "Supervisão das atividades desenvolvidas pelo aluno durante o estágio."
$ws.Range("C19").Value = "Supervisão das atividades desenvolvidas pelo aluno durante o estágio."

$ws.Range("A20").Value = "Critério:"
$ws.Range("B20").Value = "MF = Nota baseada em relatório final e no desempenho no estágio, a ser atribuída pelo professor orientador do estágio."
$ws.Range("C20").Value = "MF = Nota baseada em relatório final e no desempenho no estágio, a ser atribuída pelo professor orientador do estágio."

$ws.Range("A21").Value = "Norma de recuperação:"
$ws.Range("B21").Value = "Não será oferecida recuperação."
$ws.Range("C21").Value = "Não será oferecida recuperação."

$ws.Range("A22").Value = "Bibliografia:"
$ws.Range("B22").Value = "A ser definida com o orientador em função das atividades desenvolvidas no estágio."
$ws.Range("C22").Value = "A ser definida com o orientador em função das atividades desenvolvidas no estágio."

$ws.Range("A23").Value = "Requisitos:"

$ws.Range("B24").Value = "LOB1003 -  Cálculo I  (Requisito)`n"
$ws.Range("C24").Value = "LOB1003 -  Cálculo I  (Requisito)`n"

$ws.Range("B25").Value = "LOB1004 -  Cálculo II  (Requisito)`n"
$ws.Range("C25").Value = "LOB1004 -  Cálculo II  (Requisito)`n"

$ws.Range("B26").Value = "LOB1006 -  Cálculo IV  (Requisito)`n"
$ws.Range("C26").Value = "LOB1006 -  Cálculo IV  (Requisito)`n"

$ws.Range("B27").Value = "LOB1009 -  Leitura e Interpretação de Desenho Técnico  (Requisito)`n"
$ws.Range("C27").Value = "LOB1009 -  Leitura e Interpretação de Desenho Técnico  (Requisito)`n"

$ws.Range("B28").Value = "LOB1012 -  Estatística  (Requisito)`n"
$ws.Range("C28").Value = "LOB1012 -  Estatística  (Requisito)`n"

$ws.Range("B29").Value = "LOB1018 -  Física I  (Requisito)`n"
$ws.Range("C29").Value = "LOB1018 -  Física I  (Requisito)`n"

$ws.Range("B30").Value = "LOB1019 -  Física II  (Requisito)`n"
$ws.Range("C30").Value = "LOB1019 -  Física II  (Requisito)`n"

$ws.Range("B31").Value = "LOB1024 -  Mecânica  (Requisito)`n"
$ws.Range("C31").Value = "LOB1024 -  Mecânica  (Requisito)`n"

$ws.Range("B32").Value = "LOB1036 -  Geometria Analítica  (Requisito)`n"
$ws.Range("C32").Value = "LOB1036 -  Geometria Analítica  (Requisito)`n"

$ws.Range("B33").Value = "LOB1037 -  Àlgebra Linear  (Requisito)`n"
$ws.Range("C33").Value = "LOB1037 -  Àlgebra Linear  (Requisito)`n"

$ws.Range("B34").Value = "LOB1038 -  Física Experimental I  (Requisito)`n"
$ws.Range("C34").Value = "LOB1038 -  Física Experimental I  (Requisito)`n"

$ws.Range("B35").Value = "LOB1039 -  Física Experimental III  (Requisito)`n"
$ws.Range("C35").Value = "LOB1039 -  Física Experimental III  (Requisito)`n"

$ws.Range("B36").Value = "LOB1041 -  Física Experimental II  (Requisito)`n"
$ws.Range("C36").Value = "LOB1041 -  Física Experimental II  (Requisito)`n"

$ws.Range("B37").Value = "LOB1045 -  Leitura e Produção de Textos Acadêmicos  (Requisito)`n"
$ws.Range("C37").Value = "LOB1045 -  Leitura e Produção de Textos Acadêmicos  (Requisito)`n"

$ws.Range("B38").Value = "LOB1052 -  Cálculo III  (Requisito)`n"
$ws.Range("C38").Value = "LOB1052 -  Cálculo III  (Requisito)`n"

$ws.Range("B39").Value = "LOB1053 -  Física III  (Requisito)`n"
$ws.Range("C39").Value = "LOB1053 -  Física III  (Requisito)`n"

$ws.Range("B40").Value = "LOB1056 -  Introdução aos Métodos Numéricos e Computacionais  (Requisito)`n"
$ws.Range("C40").Value = "LOB1056 -  Introdução aos Métodos Numéricos e Computacionais  (Requisito)`n"

$ws.Range("B41").Value = "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)`n"
$ws.Range("C41").Value = "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)`n"

$ws.Range("B42").Value = "LOQ4251 -  Fundamentos de Química  (Requisito)`n"
$ws.Range("C42").Value = "LOQ4251 -  Fundamentos de Química  (Requisito)`n"

# Step 3: remove cells that no longer hold content in the new layout
$ws.Range("A13").Clear() | Out-Null
$ws.Range("B18").Clear() | Out-Null
$ws.Range("C18").Clear() | Out-Null
$ws.Range("B23").Clear() | Out-Null
$ws.Range("C23").Clear() | Out-Null

# Step 4: fix up row heights to match the new layout
$ws.Rows.Item(13).AutoFit() | Out-Null
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(17).RowHeight = 120
$ws.Rows.Item(18).AutoFit() | Out-Null
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 120
$ws.Rows.Item(23).AutoFit() | Out-Null
$ws.Rows.Item(42).RowHeight = 30
